$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.144.01"
$ws.Range("E2").Value = "  -2.13%  "

$ws.Range("D3").Value = "3.069.31"
$ws.Range("E3").Value = "  -2.26%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'520.64"
$ws.Range("E5").Value = "  -2.28%  "

$ws.Range("D6").Value = "'135.43"
$ws.Range("E6").Value = "  -5.30%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "3.069.37"
$ws.Range("E8").Value = "  -2.28%  "

$ws.Range("D9").Value = "'0.471"
$ws.Range("E9").Value = "  +5.65%  "

$ws.Range("E10").Value = "  +1.39%  "

$ws.Range("E11").Value = "  -3.35%  "

$ws.Range("D12").Value = "'0.400"
$ws.Range("E12").Value = "  +1.74%  "

$ws.Range("E13").Value = "  +1.40%  "

$ws.Range("D14").Value = "3.595.46"
$ws.Range("E14").Value = "  -2.36%  "

$ws.Range("D15").Value = "'25.01"
$ws.Range("E15").Value = "  -2.50%  "

$ws.Range("E16").Value = "  -4.44%  "

$ws.Range("D17").Value = "57.194.56"
$ws.Range("E17").Value = "  -2.10%  "

$ws.Range("D18").Value = "3.062.80"
$ws.Range("E18").Value = "  -2.58%  "

$ws.Range("E19").Value = "  -4.59%  "

$ws.Range("D20").Value = "'12.38"
$ws.Range("E20").Value = "  -3.79%  "

$ws.Range("E21").Value = "  -3.02%  "

$ws.Range("D22").Value = "'347.87"
$ws.Range("E22").Value = "  +1.05%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").Value = "'68.81"
$ws.Range("E24").Value = "  +1.64%  "

$ws.Range("D25").Value = "'0.497"
$ws.Range("E25").Value = "  -3.06%  "

$ws.Range("E26").Value = "  +0.34%  "

$ws.Range("E27").Value = "  -3.22%  "

$ws.Range("D28").Value = "0.0₃0841"
$ws.Range("E28").Value = "  -10.17%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("D30").Value = "'7.10"
$ws.Range("E30").Value = "  -5.19%  "

$ws.Range("E31").Value = "  -3.04%  "

$ws.Range("D32").Value = "'20.83"
$ws.Range("E32").Value = "  -1.36%  "

$ws.Range("E33").Value = "  -10.69%  "

$ws.Range("D34").Value = "'158.34"
$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("D35").Value = "'4.78"
$ws.Range("E35").Value = "  -0.46%  "

$ws.Range("D36").Value = "'1.12"
$ws.Range("E36").Value = "  -6.84%  "

$ws.Range("D37").Value = "'5.96"
$ws.Range("E37").Value = "  -4.39%  "

$ws.Range("D38").Value = "'25.14"
$ws.Range("E38").Value = "  -4.64%  "

$ws.Range("E39").Value = "  -3.14%  "

$ws.Range("D40").Value = "'0.0652"
$ws.Range("E40").Value = "  -2.81%  "

$ws.Range("E41").Value = "  -6.15%  "

$ws.Range("D42").Value = "'4.01"
$ws.Range("E42").Value = "  -0.14%  "

$ws.Range("D44").Value = "2.407.76"
$ws.Range("E44").Value = "  +5.34%  "

$ws.Range("D45").Value = "'36.47"
$ws.Range("E45").Value = "  -0.12%  "

$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").Value = "3.108.59"
$ws.Range("E47").Value = "  -2.29%  "

$ws.Range("D48").Value = "'0.0259"
$ws.Range("E48").Value = "  -1.78%  "

$ws.Range("D49").Value = "'5.95"
$ws.Range("E49").Value = "  -2.55%  "

$ws.Range("E50").Value = "  -8.56%  "

$ws.Range("D51").Value = "'19.25"
$ws.Range("E51").Value = "  -7.15%  "
